$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$items = @(
    @("Dwemer Kragen Helm",      "Heavy Helm",      "Exotic"),
    @("Dwemer Kragen Cuirass",   "Heavy Cuirass",   "Exotic"),
    @("Dwemer Kragen Gauntlets", "Heavy Gauntlets", "Exotic"),
    @("Dwemer Kragen Boots",     "Heavy Boots",     "Exotic"),
    @("Dwemer Kragen Shield",    "Heavy Shield",    "Exotic")
)

$row = 48
foreach ($item in $items) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $row = $row + 1
}

$ws.Range("G50").Select()
